$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & 1h volume change), and row 49/50 swap
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.649.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.54%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.904.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.11%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.29"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.68%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.59%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.913.89"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.83%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.19%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.412.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.11%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.629.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.70%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.913.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000141"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.84%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.04"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.17%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.59"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.87%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.37%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.69"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.90"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.455"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.181"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0859"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.86%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.68"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.27%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.27"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.69%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.59"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.87%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.47%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.09%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.73%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.73"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.298.88"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.88%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0582"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.55"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.73%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0237"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.16%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.32"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.39%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0926"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.87%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "251.94"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.01%  "
